$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number, Report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/6/2023  Through  2/12/2023"

# --- Crime stat table updates (rows 15-30) ---

# Row 15
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("L15").Copy($ws.Range("N15"))
$ws.Range("N15").Value = 0

# Row 16
$ws.Range("F15").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -75
$ws.Range("I16").Value = 8
$ws.Range("J16").Value = 14
$ws.Range("K16").Value = -42.857142857142
$ws.Range("N16").Value = -92.233009708737

# Row 17
$ws.Range("F15").Copy($ws.Range("C17"))
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -83.333333333333
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -42.857142857142
$ws.Range("I17").Value = 20
$ws.Range("J17").Value = 26
$ws.Range("K17").Value = -23.076923076923
$ws.Range("L17").Value = 53.846153846153
$ws.Range("M17").Value = 81.818181818181
$ws.Range("N17").Value = -33.333333333333

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("C14").Copy($ws.Range("D18"))
$ws.Range("E14").Copy($ws.Range("E18"))
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 75
$ws.Range("I18").Value = 16
$ws.Range("K18").Value = 6.666666666666
$ws.Range("L18").Value = 23.076923076923
$ws.Range("M18").Value = -71.929824561403
$ws.Range("N18").Value = -93.248945147679

# Row 19
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 27.272727272727
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = 30
$ws.Range("I19").Value = 73
$ws.Range("J19").Value = 59
$ws.Range("K19").Value = 23.728813559322
$ws.Range("L19").Value = 114.705882352941
$ws.Range("M19").Value = 52.083333333333
$ws.Range("N19").Value = -12.048192771084

# Row 20
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 66.666666666666
$ws.Range("I20").Value = 16
$ws.Range("J20").Value = 10
$ws.Range("K20").Value = 60
$ws.Range("L20").Value = 60
$ws.Range("N20").Value = -93.043478260869

# Row 21
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -4.347826086956
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 90
$ws.Range("H21").Value = 7.777777777777
$ws.Range("I21").Value = 135
$ws.Range("J21").Value = 125
$ws.Range("K21").Value = 8
$ws.Range("L21").Value = 70.886075949367
$ws.Range("M21").Value = -14.012738853503
$ws.Range("N21").Value = -80.320699708454

# Row 24
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = -10
$ws.Range("F24").Value = 70
$ws.Range("G24").Value = 68
$ws.Range("H24").Value = 2.941176470588
$ws.Range("I24").Value = 125
$ws.Range("J24").Value = 109
$ws.Range("K24").Value = 14.678899082568
$ws.Range("L24").Value = 50.602409638554
$ws.Range("M24").Value = 0.806451612903

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = -11.111111111111
$ws.Range("I25").Value = 33
$ws.Range("J25").Value = 35
$ws.Range("K25").Value = -5.714285714285
$ws.Range("L25").Value = -8.333333333333
$ws.Range("M25").Value = -25

# Row 26
$ws.Range("C14").Copy($ws.Range("C26"))
$ws.Range("F15").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("L15").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100
$ws.Range("F15").Copy($ws.Range("G26"))
$ws.Range("G26").Value = 1
$ws.Range("L15").Copy($ws.Range("H26"))
$ws.Range("H26").Value = 100
$ws.Range("J26").Value = 3
$ws.Range("K26").Value = -33.333333333333

# Row 27
$ws.Range("F15").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = -66.666666666666
$ws.Range("L27").Value = -50

# Row 30
$ws.Range("G30").Value = 2
$ws.Range("J30").Value = 2
$ws.Range("K30").Value = -50
